# The commit swaps the content of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
#   - theme1.xml goes from being empty/orphaned to holding the "Integral" /
#     "Red Violet" theme (12 colors: 000000, FFFFFF, 454551, D8D9DC, E32D91,
#     C830CC, 4EA6DC, 4775E7, 8971E1, D54773, 6B9F25, 8C8C8C).
#   - theme2.xml (the theme actually used by the single slide master, and
#     therefore by every slide in the deck) goes from "Integral" / "Red
#     Violet" to the plain built-in "Office Theme" / "Office" color scheme
#     (12 colors: 000000, FFFFFF, 44546A, E7E6E6, 5B9BD5, ED7D31, A5A5A5,
#     FFC000, 4472C4, 70AD47, 0563C1, 954F72).
#
# The font scheme and format scheme (fills/lines/effects) are byte-for-byte
# identical between the two themes already, so the only externally visible
# effect of the swap is this 12-entry colour-scheme change on the theme that
# backs the presentation's slide master. We reproduce that via the standard
# ThemeColorScheme object, which is the supported PowerPoint automation
# surface for editing a theme's colours.

$p = $ppt.ActivePresentation

# RGB() isn't available outside VBA, so pass the packed 0xBBGGRR-style
# integers directly (same encoding VBA's RGB(r,g,b) produces: r + g*256 + b*65536).
$officeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
